# BGU-14 RegNLic: new TypeEditors rows for the Appx4 physical-person questionnaire
# Adds 19 new ClassName entries (rows 33-51) to the "TypeEditors" sheet, following
# the exact same generator-formula pattern used by the existing rows (1-32).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TypeEditors")

$newClasses = @(
    "BreachOfLawRecordInfo",
    "EducationRecordInfo",
    "EmploymentRecordInfo",
    "FinancialGuaranteeInfo",
    "IncomeOriginInfo",
    "IndebtnessInfo",
    "IndebtnessInfoBase",
    "LiquidatedEntityOwnershipInfo",
    "LoanInfo",
    "PaymentDeadlineInfo",
    "PaymentModeInfo",
    "ProfessionLicenseInfo",
    "SharesAcquisitionInfo",
    "BankAccountInfo",
    "ProfessionLicensingBodyInfo",
    "PublicationInfo",
    "PublishingHouseInfo",
    "UniversityOrCollegeInfo",
    "FinancialOversightAuthorityInfo"
)

$startRow = 33
for ($i = 0; $i -lt $newClasses.Count; $i++) {
    $r = $startRow + $i
    $className = $newClasses[$i]

    $ws.Range("A$r").Value = $className

    $ws.Range("B$r").Formula = '="I"& TRIM(A' + $r + ') & "EditFormFactory"'
    $ws.Range("C$r").Formula = '="public interface I"& TRIM(A' + $r + ') & "EditFormFactory : ITypeEditorFormFactoryBase { }"'
    $ws.Range("D$r").Formula = '=A' + $r + '& "_Editor"'
    $ws.Range("E$r").Formula = '="public class " & D' + $r + ' & " : GenericTypeEditor<"&A' + $r + '&"> { private " & B' + $r + ' & " _fact; protected override ITypeEditorFormFactoryBase TypeEditorFormFactory { get { if (_fact == null) _fact = TypeEditorsDispatcher.Container.Resolve<" &B' + $r + ' & ">(); return _fact; } }  }"'
    $ws.Range("F$r").Formula = '=A' + $r + '& "EditFormFactoryBasic"'
    $ws.Range("G$r").Formula = '="public class " &F' + $r + '& " : " & B' + $r + ' & " { public System.Windows.Forms.Form SpawnInstance() { return new DummyForm<" &A' + $r + '&" >(); } }"'
    $ws.Range("H$r").Formula = '="cont.RegisterInstance<" & B' + $r + ' & ">(new " & F' + $r + '& "(), new ContainerControlledLifetimeManager());"'
    $ws.Range("I$r").Formula = '="[System.ComponentModel.Editor(typeof(BGU.DRPL.SignificantOwnership.Core.TypeEditors." &D' + $r + ' & "), typeof(System.Drawing.Design.UITypeEditor))]"'
}

$lastRow = $startRow + $newClasses.Count - 1

$ws.Activate()
$ws.Range("G$lastRow").Select()
